# MONTENEGRO.xlsx — sheet cleanup:
#   - "Paineis DARQ"            -> "PAINEIS DARQ"
#   - "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#   - remove the now-unneeded "Desarquivamentos Pendentes" sheet

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$null = $wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Deleting the sheet shifts Excel's active-sheet cursor; restore it so the
# originally-selected tab ("PAINEIS DARQ") stays tabSelected="1".
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
